$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

$ws.Range("I2:I10").NumberFormat = "@"

$ws.Cells.Item(1, 2).Value = 'bank'
$ws.Cells.Item(1, 3).Value = 'deposit_type'
$ws.Cells.Item(1, 4).Value = 'currency'
$ws.Cells.Item(1, 5).Value = 'owner'
$ws.Cells.Item(1, 6).Value = 'total'
$ws.Cells.Item(1, 7).Value = 'property_category'
$ws.Cells.Item(1, 8).Value = 'category'
$ws.Cells.Item(1, 9).Value = 'date'
$ws.Cells.Item(1, 10).Value = 'legislator_name'
$ws.Cells.Item(1, 11).Value = 'legislator_id'
$ws.Cells.Item(1, 12).Value = 'source_file'
$ws.Cells.Item(1, 13).Value = 'index'
$ws.Cells.Item(2, 1).Value = 75
$ws.Cells.Item(2, 2).Value = '聯邦商業銀行營業部'
$ws.Cells.Item(2, 3).Value = '活期儲蓄存款'
$ws.Cells.Item(2, 4).Value = '新臺幣'
$ws.Cells.Item(2, 5).Value = '陳怡潔'
$ws.Cells.Item(2, 6).Value = 20220316
$ws.Cells.Item(2, 7).Value = 'deposit'
$ws.Cells.Item(2, 8).Value = 'normal'
$ws.Cells.Item(2, 9).Value = '2013-06-20'
$ws.Cells.Item(2, 10).Value = '陳怡潔'
$ws.Cells.Item(2, 11).Value = 1804
$ws.Cells.Item(2, 12).Value = 'tmp20f31'
$ws.Cells.Item(2, 13).Value = 75
$ws.Cells.Item(3, 1).Value = 76
$ws.Cells.Item(3, 2).Value = '臺灣土地銀行新莊分行'
$ws.Cells.Item(3, 3).Value = '活期儲蓄存款'
$ws.Cells.Item(3, 4).Value = '新臺幣'
$ws.Cells.Item(3, 5).Value = '陳怡潔'
$ws.Cells.Item(3, 6).Value = 109162
$ws.Cells.Item(3, 7).Value = 'deposit'
$ws.Cells.Item(3, 8).Value = 'normal'
$ws.Cells.Item(3, 9).Value = '2013-06-20'
$ws.Cells.Item(3, 10).Value = '陳怡潔'
$ws.Cells.Item(3, 11).Value = 1804
$ws.Cells.Item(3, 12).Value = 'tmp20f31'
$ws.Cells.Item(3, 13).Value = 76
$ws.Cells.Item(4, 1).Value = 77
$ws.Cells.Item(4, 2).Value = '臺灣銀行淡水分行'
$ws.Cells.Item(4, 3).Value = '活期儲蓄存款'
$ws.Cells.Item(4, 4).Value = '新臺幣'
$ws.Cells.Item(4, 5).Value = '陳怡潔'
$ws.Cells.Item(4, 6).Value = 663074
$ws.Cells.Item(4, 7).Value = 'deposit'
$ws.Cells.Item(4, 8).Value = 'normal'
$ws.Cells.Item(4, 9).Value = '2013-06-20'
$ws.Cells.Item(4, 10).Value = '陳怡潔'
$ws.Cells.Item(4, 11).Value = 1804
$ws.Cells.Item(4, 12).Value = 'tmp20f31'
$ws.Cells.Item(4, 13).Value = 77
$ws.Cells.Item(5, 1).Value = 78
$ws.Cells.Item(5, 2).Value = '合作金庫商業銀行南西分行'
$ws.Cells.Item(5, 3).Value = '活期儲蓄存款'
$ws.Cells.Item(5, 4).Value = '新臺幣'
$ws.Cells.Item(5, 5).Value = '陳怡潔'
$ws.Cells.Item(5, 6).Value = 96980
$ws.Cells.Item(5, 7).Value = 'deposit'
$ws.Cells.Item(5, 8).Value = 'normal'
$ws.Cells.Item(5, 9).Value = '2013-06-20'
$ws.Cells.Item(5, 10).Value = '陳怡潔'
$ws.Cells.Item(5, 11).Value = 1804
$ws.Cells.Item(5, 12).Value = 'tmp20f31'
$ws.Cells.Item(5, 13).Value = 78
$ws.Cells.Item(6, 1).Value = 79
$ws.Cells.Item(6, 2).Value = '臺灣銀行群賢分行'
$ws.Cells.Item(6, 3).Value = '定期儲蓄存款'
$ws.Cells.Item(6, 4).Value = '新臺幣'
$ws.Cells.Item(6, 5).Value = '陳怡潔'
$ws.Cells.Item(6, 6).Value = 406745
$ws.Cells.Item(6, 7).Value = 'deposit'
$ws.Cells.Item(6, 8).Value = 'normal'
$ws.Cells.Item(6, 9).Value = '2013-06-20'
$ws.Cells.Item(6, 10).Value = '陳怡潔'
$ws.Cells.Item(6, 11).Value = 1804
$ws.Cells.Item(6, 12).Value = 'tmp20f31'
$ws.Cells.Item(6, 13).Value = 79
$ws.Cells.Item(7, 1).Value = 80
$ws.Cells.Item(7, 2).Value = '國泰世華商業銀行建成分行'
$ws.Cells.Item(7, 3).Value = '綜合存款'
$ws.Cells.Item(7, 4).Value = '新臺幣'
$ws.Cells.Item(7, 5).Value = '陳怡潔'
$ws.Cells.Item(7, 6).Value = 593732
$ws.Cells.Item(7, 7).Value = 'deposit'
$ws.Cells.Item(7, 8).Value = 'normal'
$ws.Cells.Item(7, 9).Value = '2013-06-20'
$ws.Cells.Item(7, 10).Value = '陳怡潔'
$ws.Cells.Item(7, 11).Value = 1804
$ws.Cells.Item(7, 12).Value = 'tmp20f31'
$ws.Cells.Item(7, 13).Value = 80
$ws.Cells.Item(8, 1).Value = 81
$ws.Cells.Item(8, 2).Value = '臺灣土地銀行國外部'
$ws.Cells.Item(8, 3).Value = '綜合存款'
$ws.Cells.Item(8, 4).Value = '美金'
$ws.Cells.Item(8, 5).Value = '陳怡潔'
$ws.Cells.Item(8, 6).Value = 3000
$ws.Cells.Item(8, 7).Value = 'deposit'
$ws.Cells.Item(8, 8).Value = 'normal'
$ws.Cells.Item(8, 9).Value = '2013-06-20'
$ws.Cells.Item(8, 10).Value = '陳怡潔'
$ws.Cells.Item(8, 11).Value = 1804
$ws.Cells.Item(8, 12).Value = 'tmp20f31'
$ws.Cells.Item(8, 13).Value = 81
$ws.Cells.Item(9, 1).Value = 82
$ws.Cells.Item(9, 2).Value = '國泰世華商業銀行國外部'
$ws.Cells.Item(9, 3).Value = '活期存款'
$ws.Cells.Item(9, 4).Value = '美金'
$ws.Cells.Item(9, 5).Value = '陳怡潔'
$ws.Cells.Item(9, 6).Value = 123515
$ws.Cells.Item(9, 7).Value = 'deposit'
$ws.Cells.Item(9, 8).Value = 'normal'
$ws.Cells.Item(9, 9).Value = '2013-06-20'
$ws.Cells.Item(9, 10).Value = '陳怡潔'
$ws.Cells.Item(9, 11).Value = 1804
$ws.Cells.Item(9, 12).Value = 'tmp20f31'
$ws.Cells.Item(9, 13).Value = 82
$ws.Cells.Item(10, 1).Value = 83
$ws.Cells.Item(10, 2).Value = '曰盛國際商業銀行延平分行'
$ws.Cells.Item(10, 3).Value = '綜合存款'
$ws.Cells.Item(10, 4).Value = '新臺幣'
$ws.Cells.Item(10, 5).Value = '陳◦霓'
$ws.Cells.Item(10, 6).Value = 2445408
$ws.Cells.Item(10, 7).Value = 'deposit'
$ws.Cells.Item(10, 8).Value = 'normal'
$ws.Cells.Item(10, 9).Value = '2013-06-20'
$ws.Cells.Item(10, 10).Value = '陳怡潔'
$ws.Cells.Item(10, 11).Value = 1804
$ws.Cells.Item(10, 12).Value = 'tmp20f31'
$ws.Cells.Item(10, 13).Value = 83
